$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the "Tableau4" table by two rows (ListRows.Add keeps the table's
# styling/dxf wiring and grows ref="A1:H94" -> "A1:H96" like Excel does).
$lo = $ws.ListObjects.Item("Tableau4")
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# Copy the formatting of the last "real" data row (94) down onto the two
# freshly added rows so the new cells pick up styles s=21/24/24/20/5... like
# row 94, instead of the generic column default styles.
$ws.Range("A94:H94").Copy() | Out-Null
$ws.Range("A95:H96").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 95: 2022-04-01, 17:40 -> 18:30, CPNV / Diagramme de flux / "J'ai commencé..."
$ws.Range("A95").Value = 44652
$ws.Range("B95").Value = 0.73611111111111116
$ws.Range("C95").Value = 0.77083333333333337
$ws.Range("D95").Formula = "=Tableau4[[#This Row],[Heure fin]]-Tableau4[[#This Row],[Heure début]]"
$ws.Range("E95").Value = "CPNV"
$ws.Range("F95").Value = "Diagramme de flux"
$ws.Range("G95").Value = "J'ai commencé le diagramme de flux"

# Row 96: 2022-04-01, 14:10 -> 14:59, CPNV / Diagramme de flux / "J'ai continué..."
$ws.Range("A96").Value = 44652
$ws.Range("B96").Value = 0.59027777777777779
$ws.Range("C96").Value = 0.62430555555555556
$ws.Range("D96").Formula = "=Tableau4[[#This Row],[Heure fin]]-Tableau4[[#This Row],[Heure début]]"
$ws.Range("E96").Value = "CPNV"
$ws.Range("F96").Value = "Diagramme de flux"
$ws.Range("G96").Value = "J'ai continué le diagramme de flux"

# Scroll/selection bookkeeping to mirror the author's final cursor position.
$ws.Range("G107").Select()
$excel.ActiveWindow.ScrollRow = 70
$excel.ActiveWindow.ScrollColumn = 1
